# Update SmokeTests to add conda and Personal Access Token
# - Move the "Assignee" column from the "issues" sheet to the "OS instructions" sheet
# - Remove the "Conda Package Tests" row from the "issues" sheet
# - Update the active sheet/selection to "OS instructions"

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("OS instructions")
$ws2 = $wb.Worksheets.Item("issues")

# --- "OS instructions" sheet: add a new "Assignee" header in column C ---
$c1 = $ws1.Range("C1")
$c1.Value = "Assignee"
$c1.Font.Bold = $true

# --- "issues" sheet: remove the old "Assignee" header cell (column D) ---
$ws2.Range("D1").Clear()

# --- "issues" sheet: remove the "Conda Package Tests" smoke test row (row 8) ---
$ws2.Range("A8:B8").ClearContents()
$ws2.Range("C8").Clear()
$ws2.Rows.Item(8).RowHeight = 16

# --- Update selections / active sheet ---
[void]$ws2.Range("D1").Select()
[void]$ws1.Select()
[void]$ws1.Range("C1").Select()
